# Update sample Excel file with tested column additions
# - Removes the old "CONDITION 1" (age) column entirely
# - Keeps the old "CONDITION 2" (status) column, which shifts left into column B
# - Inserts a brand-new "New Credit Check" CONDITION column at C
# - Inserts a brand-new "Loyalty Reward" ACTION column at E (pushing the
#   existing "Bonus Action" column to F)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the CONDITION 1 (age) column - column B - entirely.
#    This shifts the old CONDITION 2 (status) column from C to B,
#    ACTION 1 from D to C, and Bonus Action from E to D.
$ws.Columns.Item(2).EntireColumn.Delete()

# 2. Insert a new blank column at C for the new "New Credit Check" condition.
#    This pushes ACTION 1 back to D and Bonus Action back to E.
$ws.Columns.Item(3).EntireColumn.Insert()

# 3. Insert a new blank column at E for the new "Loyalty Reward" action.
#    This pushes Bonus Action from E to F.
$ws.Columns.Item(5).EntireColumn.Insert()

# 4. Populate the new "New Credit Check" condition column (C)
$ws.Range("C1").Value = "New Credit Check"
$ws.Range("C2").Value = "customer.getCreditScore() >= `$param"

# 5. Populate the new "Loyalty Reward" action column (E)
$ws.Range("E1").Value = "Loyalty Reward"
$ws.Range("E2").Value = "customer.setLoyaltyPoints(`$param);"
